# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the Sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as captured by the scheduled runner snapshot.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 936.4722
$ws.Range("I19").Value = 700.75
$ws.Range("J19").Value = 1054.3334
$ws.Range("K19").Value = 700.75
$ws.Range("L19").Value = 1054.3334
$ws.Range("M19").Value = -525.75
$ws.Range("N19").Value = -1404.3334

$ws.Range("H38").Value = 471.14816
$ws.Range("I38").Value = 155.04167
$ws.Range("K38").Value = 465.12501
$ws.Range("M38").Value = -93.12501000000003

$ws.Range("H76").Value = 3099.3076
$ws.Range("I76").Value = 3098.6316
$ws.Range("J76").Value = 3125
$ws.Range("K76").Value = 3098.6316
$ws.Range("L76").Value = 3125
$ws.Range("M76").Value = -2783.6316
$ws.Range("N76").Value = -3755

$ws.Range("H79").Value = 3099.3076
$ws.Range("I79").Value = 3098.6316
$ws.Range("J79").Value = 3125
$ws.Range("K79").Value = 3098.6316
$ws.Range("L79").Value = 3125
$ws.Range("M79").Value = -2006.6316
$ws.Range("N79").Value = -5309

$ws.Range("H98").Value = 1792
$ws.Range("I98").Value = 1792
$ws.Range("K98").Value = 1792
$ws.Range("M98").Value = -294

$ws.Range("H122").Value = 1792
$ws.Range("I122").Value = 1792
$ws.Range("K122").Value = 5376
$ws.Range("M122").Value = -2926

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3856.111
$ws.Range("I63").Value = 2568.3333
$ws.Range("J63").Value = 4500
$ws.Range("K63").Value = 2568.3333
$ws.Range("L63").Value = 4500
$ws.Range("M63").Value = -1882.3333
$ws.Range("N63").Value = -5872

$ws.Range("H66").Value = 3856.111
$ws.Range("I66").Value = 2568.3333
$ws.Range("J66").Value = 4500
$ws.Range("K66").Value = 12841.6665
$ws.Range("L66").Value = 22500
$ws.Range("M66").Value = -9409.666499999999
$ws.Range("N66").Value = -29364

$ws.Range("H74").Value = 32697.229
$ws.Range("I74").Value = 36835.32
$ws.Range("K74").Value = 36835.32
$ws.Range("M74").Value = -35961.32

$ws.Range("H77").Value = 32697.229
$ws.Range("I77").Value = 36835.32
$ws.Range("K77").Value = 184176.6
$ws.Range("M77").Value = -179808.6

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7778.96
$ws.Range("I134").Value = 8310.875
$ws.Range("K134").Value = 24932.625
$ws.Range("M134").Value = -22397.625

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2299.0967
$ws.Range("I31").Value = 1452.85
$ws.Range("J31").Value = 3837.7273
$ws.Range("K31").Value = 1452.85
$ws.Range("L31").Value = 3837.7273
$ws.Range("M31").Value = -1157.85
$ws.Range("N31").Value = -4427.7273

$ws.Range("H34").Value = 2299.0967
$ws.Range("I34").Value = 1452.85
$ws.Range("J34").Value = 3837.7273
$ws.Range("K34").Value = 1452.85
$ws.Range("L34").Value = 3837.7273
$ws.Range("M34").Value = -1250.85
$ws.Range("N34").Value = -4241.7273

$ws.Range("H87").Value = 26900
$ws.Range("J87").Value = 26900
$ws.Range("L87").Value = 26900
$ws.Range("N87").Value = -29272

$ws.Range("H90").Value = 26900
$ws.Range("J90").Value = 26900
$ws.Range("L90").Value = 80700
$ws.Range("N90").Value = -92556

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H122").Value = 864.1667
$ws.Range("I122").Value = 719.61536
$ws.Range("J122").Value = 1240
$ws.Range("K122").Value = 2158.84608
$ws.Range("L122").Value = 3720
$ws.Range("M122").Value = 291.1539199999997
$ws.Range("N122").Value = -8620

$ws.Range("H141").Value = 58243.25
$ws.Range("I141").Value = 60296
$ws.Range("J141").Value = 57950
$ws.Range("K141").Value = 60296
$ws.Range("L141").Value = 57950
$ws.Range("M141").Value = -55116
$ws.Range("N141").Value = -68310

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 17356292
$ws.Range("I114").Value = 25000352
$ws.Range("J114").Value = 12988257
$ws.Range("K114").Value = 75001056
$ws.Range("L114").Value = 38964771
$ws.Range("M114").Value = -74997802
$ws.Range("N114").Value = -38971279

$ws.Range("H117").Value = 744
$ws.Range("I117").Value = 400
$ws.Range("J117").Value = 916
$ws.Range("K117").Value = 1200
$ws.Range("L117").Value = 2748
$ws.Range("M117").Value = 2242
$ws.Range("N117").Value = -9632

$ws.Range("H121").Value = 1573.1818
$ws.Range("J121").Value = 1573.1818
$ws.Range("L121").Value = 4719.5454
$ws.Range("N121").Value = -7339.5454

$ws.Range("H129").Value = 1926.0555
$ws.Range("I129").Value = 1488.8667
$ws.Range("J129").Value = 2238.3333
$ws.Range("K129").Value = 4466.6001
$ws.Range("L129").Value = 6714.999899999999
$ws.Range("M129").Value = 533.3999000000003
$ws.Range("N129").Value = -16714.9999

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2149.4736
$ws.Range("I102").Value = 2210.7693
$ws.Range("K102").Value = 2210.7693
$ws.Range("M102").Value = -588.7692999999999

$ws.Range("H126").Value = 1924.3077
$ws.Range("I126").Value = 750.8333
$ws.Range("J126").Value = 2930.1428
$ws.Range("K126").Value = 2252.4999
$ws.Range("L126").Value = 8790.428400000001
$ws.Range("M126").Value = 217.5001000000002
$ws.Range("N126").Value = -13730.4284

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1183
$ws.Range("I82").Value = 963
$ws.Range("J82").Value = 1263
$ws.Range("K82").Value = 963
$ws.Range("L82").Value = 1263
$ws.Range("M82").Value = -602
$ws.Range("N82").Value = -1985

$ws.Range("H85").Value = 1183
$ws.Range("I85").Value = 963
$ws.Range("J85").Value = 1263
$ws.Range("K85").Value = 963
$ws.Range("L85").Value = 1263
$ws.Range("M85").Value = 285
$ws.Range("N85").Value = -3759

$ws.Range("H122").Value = 2490.282
$ws.Range("I122").Value = 1931.7916
$ws.Range("J122").Value = 3383.8667
$ws.Range("K122").Value = 5795.3748
$ws.Range("L122").Value = 10151.6001
$ws.Range("M122").Value = -3345.3748
$ws.Range("N122").Value = -15051.6001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 34951.445
$ws.Range("I122").Value = 1533.091
$ws.Range("J122").Value = 87466
$ws.Range("K122").Value = 4599.272999999999
$ws.Range("L122").Value = 262398
$ws.Range("M122").Value = -2149.272999999999
$ws.Range("N122").Value = -267298

$ws.Range("H126").Value = 1476.3334
$ws.Range("I126").Value = 638.4545000000001
$ws.Range("J126").Value = 2185.3076
$ws.Range("K126").Value = 1915.3635
$ws.Range("L126").Value = 6555.9228
$ws.Range("M126").Value = 554.6364999999998
$ws.Range("N126").Value = -11495.9228
